# Adds two new columns, "I0" (column I) and "IF" (column J), to the sheet,
# mirroring the header style used by the other header cells (B1:H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they look like the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2..46: column I (I0) and column J (IF)
$data = @(
    @(9,9),
    @(1,4),
    @(8,8),
    @(2,4),
    @(10,10),
    @(3,5),
    @(11,11),
    @(7,8),
    @(6,7),
    @(9,9),
    @(10,10),
    @(4,6),
    @(9,9),
    @(9,9),
    @(11,11),
    @(7,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(5,6),
    @(8,9),
    @(5,5),
    @(2,5),
    @(8,8),
    @(5,6),
    @(6,9),
    @(8,8),
    @(1,5),
    @(9,9),
    @(1,4),
    @(8,8),
    @(6,6),
    @(6,8),
    @(9,9),
    @(8,9),
    @(4,7),
    @(4,5),
    @(1,3),
    @(7,9),
    @(6,8),
    @(8,8),
    @(5,6),
    @(5,6),
    @(3,4),
    @(1,2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
